$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2:A7 with the new combined values
$ws.Range("A2").Value = "('Hydra', ['Token Creature — Hydra', '*/*'])"
$ws.Range("A3").Value = "('Minotaur', ['Token Creature — Minotaur', '2/3'])"
$ws.Range("A4").Value = "('Snake', ['Token Enchantment Creature — Snake', 'Deathtouch', '1/1'])"
$ws.Range("A5").Value = "('Sphinx', ['Token Creature — Sphinx', 'Flying', '4/4'])"
$ws.Range("A6").Value = "('Spider', ['Token Enchantment Creature — Spider', 'Reach', '1/3'])"
$ws.Range("A7").Value = "('Zombie', ['Token Creature — Zombie', '*/*'])"

# Clear out the now-unused rows 8 through 22
$ws.Range("A8:A22").Clear()
